# Alternative branch creation and checkout
$d = $word.ActiveDocument

# --- 1. Append new runs to the end of the first paragraph ---------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# Move the end of the range to just before the paragraph mark so new text
# is inserted inside the paragraph (not after it).
$r1.End = $r1.End - 1

$spacer = $r1.InsertAfter("  ")
$spacer.Font.Color = 4211712  # wdColorAutomatic default (unused placeholder)

# Re-fetch range after mutation and insert the colored text in three runs
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1

$run1 = $r1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run1.Font.Color = 12583104  # 0xC00000 -> BGR 0x0000C0 -> decimal 12583104

$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$run2 = $r1.InsertAfter("rsion for branch alternate")
$run2.Font.Color = 12583104

$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$run3 = $r1.InsertAfter(")")
$run3.Font.Color = 12583104

# --- 2. Add an empty paragraph right before the final sectPr ------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
